# Natmi following Dr Hou advice
# Update LR-pair worksheet: add "M1" and "M2" target clusters (so each of the
# three sending clusters ECs/FAPs/sCs now has five target clusters:
# ECs, FAPs, M1, M2, sCs) and refresh all computed statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column headers (row 1) stay the same; only data rows 2-16 change.

$data = @(
  @("ECs","Vegfc","Flt4","ECs",3,1,2.663313,7.989939000000001,0.3794306644527501,0.3794306644527502,3,1,11.68075,35.04225,0.9587437520473568,0.9587437520473567,31.10949332474999,279.98543992275,0.3637767788792514,0.3637767788792514),
  @("ECs","Vegfc","Flt4","FAPs",3,1,2.663313,7.989939000000001,0.3794306644527501,0.3794306644527502,3,1,0.163766,0.491298,0.01344174212253387,0.01344174212253387,0.436160116758,3.925441050822001,0.005100209144955545,0.005100209144955545),
  @("ECs","Vegfc","Flt4","M1",3,1,2.663313,7.989939000000001,0.3794306644527501,0.3794306644527502,1,0.3333333333333333,0.046993,0.140979,0.003857136326003163,0.003857136326003163,0.125157067809,1.126413610281,0.00146351579906022,0.00146351579906022),
  @("ECs","Vegfc","Flt4","M2",3,1,2.663313,7.989939000000001,0.3794306644527501,0.3794306644527502,1,0.3333333333333333,0.07763399999999999,0.232902,0.006372117582042636,0.006372117582042635,0.206763641442,1.860872772978,0.002417776808125489,0.002417776808125489),
  @("ECs","Vegfc","Flt4","sCs",3,1,2.663313,7.989939000000001,0.3794306644527501,0.3794306644527502,3,1,0.214248,0.642744,0.01758525192206341,0.0175852519220634,0.570609483624,5.135485352616,0.006672383821357521,0.00667238382135752),
  @("FAPs","Vegfc","Flt4","ECs",3,1,3.178631333333334,9.535894000000001,0.4528458348143826,0.4528458348143827,3,1,11.68075,35.04225,0.9587437520473568,0.9587437520473567,37.12879794683333,334.1591815215,0.4341631147689587,0.4341631147689587),
  @("FAPs","Vegfc","Flt4","FAPs",3,1,3.178631333333334,9.535894000000001,0.4528458348143826,0.4528458348143827,3,1,0.163766,0.491298,0.01344174212253387,0.01344174212253387,0.5205517389346667,4.684965650412001,0.006087036932838501,0.006087036932838501),
  @("FAPs","Vegfc","Flt4","M1",3,1,3.178631333333334,9.535894000000001,0.4528458348143826,0.4528458348143827,1,0.3333333333333333,0.046993,0.140979,0.003857136326003163,0.003857136326003163,0.1493734222473334,1.344360800226,0.001746688119541783,0.001746688119541783),
  @("FAPs","Vegfc","Flt4","M2",3,1,3.178631333333334,9.535894000000001,0.4528458348143826,0.4528458348143827,1,0.3333333333333333,0.07763399999999999,0.232902,0.006372117582042636,0.006372117582042635,0.246769864932,2.220928784388,0.002885586905975502,0.002885586905975502),
  @("FAPs","Vegfc","Flt4","sCs",3,1,3.178631333333334,9.535894000000001,0.4528458348143826,0.4528458348143827,3,1,0.214248,0.642744,0.01758525192206341,0.0175852519220634,0.681015405904,6.129138653136001,0.00796340808706803,0.00796340808706803),
  @("sCs","Vegfc","Flt4","ECs",3,1,1.177290666666667,3.531872,0.1677235007328671,0.1677235007328671,3,1,11.68075,35.04225,0.9587437520473568,0.9587437520473567,13.75163795466666,123.764741592,0.1608038583991466,0.1608038583991466),
  @("sCs","Vegfc","Flt4","FAPs",3,1,1.177290666666667,3.531872,0.1677235007328671,0.1677235007328671,3,1,0.163766,0.491298,0.01344174212253387,0.01344174212253387,0.1928001833173333,1.735201649856,0.00225449604473982,0.00225449604473982),
  @("sCs","Vegfc","Flt4","M1",3,1,1.177290666666667,3.531872,0.1677235007328671,0.1677235007328671,1,0.3333333333333333,0.046993,0.140979,0.003857136326003163,0.003857136326003163,0.05532442029866667,0.497919782688,0.00064693240740116,0.0006469324074011598),
  @("sCs","Vegfc","Flt4","M2",3,1,1.177290666666667,3.531872,0.1677235007328671,0.1677235007328671,1,0.3333333333333333,0.07763399999999999,0.232902,0.006372117582042636,0.006372117582042635,0.091397783616,0.822580052544,0.001068753867941643,0.001068753867941643),
  @("sCs","Vegfc","Flt4","sCs",3,1,1.177290666666667,3.531872,0.1677235007328671,0.1677235007328671,3,1,0.214248,0.642744,0.01758525192206341,0.0175852519220634,0.252232170752,2.270089536768,0.002949460013637855,0.002949460013637854)
)

$rowIndex = 2
foreach ($row in $data) {
  $colIndex = 1
  foreach ($val in $row) {
    $ws.Cells.Item($rowIndex, $colIndex).Value = $val
    $colIndex++
  }
  $rowIndex++
}
